$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet1 columns D (Price) and E (Volume(1h)) store literal text values
# (e.g. "303.55", "5.82%") rather than numbers, so each new value below is
# written with a leading apostrophe to force text entry - this keeps the
# cell a string instead of Excel auto-converting it to a Number/Percentage.
$quote = [char]39
$updates = @(
    @("D2", "303.12"),
    @("E2", "5.66%"),
    @("D3", "31.90"),
    @("E3", "9.28%"),
    @("D4", "5.247"),
    @("E4", "0.77%"),
    @("D5", "0.07435"),
    @("E5", "6.18%"),
    @("D6", "7.849"),
    @("E6", "5.67%"),
    @("D7", "3.806"),
    @("E7", "7.08%"),
    @("D8", "1.520"),
    @("E8", "7.48%"),
    @("D9", "0.9210"),
    @("E9", "2.14%"),
    @("E10", "2,600.97%"),
    @("D11", "0.1696"),
    @("E11", "5.55%"),
    @("D12", "0.07972"),
    @("E12", "5.81%"),
    @("D13", "0.07986"),
    @("E13", "4.05%"),
    @("D14", "0.03063"),
    @("E14", "4.49%"),
    @("D15", "0.09900"),
    @("E15", "9.81%"),
    @("D16", "0.001498"),
    @("E16", "-5.62%"),
    @("D17", "0.04611"),
    @("E17", "1.80%"),
    @("D18", "0.006235"),
    @("E18", "0.79%"),
    @("D19", "3.466"),
    @("E19", "-0.22%"),
    @("E20", "0.08%"),
    @("E21", "2.66%"),
    @("D22", "0.1329"),
    @("E22", "-0.30%"),
    @("D23", "4.492"),
    @("E23", "11.91%"),
    @("E24", "1.60%"),
    @("D25", "0.001223"),
    @("E25", "1.18%"),
    @("D26", "0.004437"),
    @("E26", "4.61%"),
    @("E27", "19.91%"),
    @("E28", "4.99%"),
    @("D40", "0.04503"),
    @("E40", "3.71%"),
    @("E41", "3.42%"),
    @("D42", "0.1347"),
    @("E42", "8.02%"),
    @("D43", "0.002233"),
    @("E43", "7.96%"),
    @("E44", "10.80%"),
    @("D45", "0.00006147"),
    @("E45", "5.19%"),
    @("E47", "-0.30%")
)

foreach ($u in $updates) {
    $ref = $u[0]
    $newValue = $u[1]
    $ws.Range($ref).Value = "$quote$newValue"
}
